$wb = $excel.ActiveWorkbook

# Top-level ObjTables header string lives on sheet 1 (Compartment), row 1 (A1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Unprotect()
$ws1.Cells.Item(1,1).Value = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.9' date='2020-04-26 21:09:42'"

# Sheet 1: Compartment
$ws1.Cells.Item(2,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compartment' name='Compartment' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 2: Compound
$ws = $wb.Worksheets.Item(2)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compound' name='Compound' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 3: Definition
$ws = $wb.Worksheets.Item(3)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Definition' name='Definition' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 4: Enzyme
$ws = $wb.Worksheets.Item(4)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Enzyme' name='Enzyme' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 5: FbcObjective
$ws = $wb.Worksheets.Item(5)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='FbcObjective' name='FbcObjective' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 6: Gene
$ws = $wb.Worksheets.Item(6)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Gene' name='Gene' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 7: Layout
$ws = $wb.Worksheets.Item(7)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Layout' name='Layout' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 8: Measurement
$ws = $wb.Worksheets.Item(8)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Measurement' name='Measurement' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 9: PbConfig
$ws = $wb.Worksheets.Item(9)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='PbConfig' name='PbConfig' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 10: Position
$ws = $wb.Worksheets.Item(10)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Position' name='Position' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 11: Protein
$ws = $wb.Worksheets.Item(11)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Protein' name='Protein' date='2020-04-26 21:09:42' objTablesVersion='0.0.9'"

# Sheet 12: Quantity
$ws = $wb.Worksheets.Item(12)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Quantity' name='Quantity' date='2020-04-26 21:09:42' objTablesVersion='0.0.9' level='1.0' version='0.1'"

# Sheet 13: QuantityInfo
$ws = $wb.Worksheets.Item(13)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityInfo' name='QuantityInfo' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 14: QuantityMatrix
$ws = $wb.Worksheets.Item(14)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityMatrix' name='QuantityMatrix' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 15: Reaction
$ws = $wb.Worksheets.Item(15)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Reaction' name='Reaction' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 16: ReactionStoichiometry
$ws = $wb.Worksheets.Item(16)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 17: Regulator
$ws = $wb.Worksheets.Item(17)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Regulator' name='Regulator' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 18: Relation
$ws = $wb.Worksheets.Item(18)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relation' name='Relation' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 19: Relationship
$ws = $wb.Worksheets.Item(19)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relationship' name='Relationship' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 20: SparseMatrix
$ws = $wb.Worksheets.Item(20)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrix' name='SparseMatrix' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 21: SparseMatrixColumn
$ws = $wb.Worksheets.Item(21)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 22: SparseMatrixOrdered
$ws = $wb.Worksheets.Item(22)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 23: SparseMatrixRow
$ws = $wb.Worksheets.Item(23)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixRow' name='SparseMatrixRow' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 24: StoichiometricMatrix
$ws = $wb.Worksheets.Item(24)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 25: rxnconContingencyList
$ws = $wb.Worksheets.Item(25)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconContingencyList' name='rxnconContingencyList' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# Sheet 26: rxnconReactionList
$ws = $wb.Worksheets.Item(26)
$ws.Unprotect()
$ws.Cells.Item(1,1).Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconReactionList' name='rxnconReactionList' date='2020-04-26 21:09:43' objTablesVersion='0.0.9'"

# --- Relation sheet (18): G2 "!From" -> "!FromObject", H2 "!To" -> "!ToObject" ---
$wsRelation = $wb.Worksheets.Item(18)
$wsRelation.Unprotect()
$wsRelation.Cells.Item(2,7).Value = "!FromObject"
$wsRelation.Cells.Item(2,8).Value = "!ToObject"

$valG = $wsRelation.Range("G2").Validation
$valG.ErrorTitle = "FromObject"
$valG.InputTitle = "FromObject"

$valH = $wsRelation.Range("H2").Validation
$valH.ErrorTitle = "ToObject"
$valH.InputTitle = "ToObject"

# --- Relationship sheet (19): B2 "!From" -> "!FromObject", C2 "!To" -> "!ToObject" ---
$wsRelationship = $wb.Worksheets.Item(19)
$wsRelationship.Unprotect()
$wsRelationship.Cells.Item(2,2).Value = "!FromObject"
$wsRelationship.Cells.Item(2,3).Value = "!ToObject"

$valB = $wsRelationship.Range("B2").Validation
$valB.ErrorTitle = "FromObject"
$valB.InputTitle = "FromObject"

$valC = $wsRelationship.Range("C2").Validation
$valC.ErrorTitle = "ToObject"
$valC.InputTitle = "ToObject"
